$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7887.636
$ws.Range("I62").Value = 7753
$ws.Range("K62").Value = 7753
$ws.Range("M62").Value = -7129
$ws.Range("H65").Value = 7887.636
$ws.Range("I65").Value = 7753
$ws.Range("K65").Value = 38765
$ws.Range("M65").Value = -35645
$ws.Range("H92").Value = 1369.5358
$ws.Range("I92").Value = 398.2857
$ws.Range("K92").Value = 398.2857
$ws.Range("M92").Value = 849.7143
$ws.Range("H99").Value = 1361.7
$ws.Range("J99").Value = 2093.6667
$ws.Range("L99").Value = 6281.000100000001
$ws.Range("N99").Value = -9277.000100000001
$ws.Range("H106").Value = 1549.6
$ws.Range("I106").Value = 1322.7693
$ws.Range("J106").Value = 3024
$ws.Range("K106").Value = 1322.7693
$ws.Range("L106").Value = 3024
$ws.Range("M106").Value = -691.7692999999999
$ws.Range("N106").Value = -4286
$ws.Range("H135").Value = 485.83334
$ws.Range("I135").Value = 485.83334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4372.50006
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1837.50006
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 34649.535
$ws.Range("I137").Value = 40908.29
$ws.Range("K137").Value = 122724.87
$ws.Range("M137").Value = -120174.87
$ws.Range("H138").Value = 2922.8823
$ws.Range("J138").Value = 3580.0852
$ws.Range("L138").Value = 10740.2556
$ws.Range("N138").Value = -21020.2556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10226.824
$ws.Range("I32").Value = 5909.2324
$ws.Range("J32").Value = 23488
$ws.Range("K32").Value = 5909.2324
$ws.Range("L32").Value = 23488
$ws.Range("M32").Value = -5622.2324
$ws.Range("N32").Value = -24062
$ws.Range("H61").Value = 4520.3687
$ws.Range("I61").Value = 4494
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 4494
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -4282
$ws.Range("N61").Value = -5419
$ws.Range("H132").Value = 2021.45
$ws.Range("I132").Value = 1842.1538
$ws.Range("K132").Value = 5526.4614
$ws.Range("M132").Value = -2996.4614
$ws.Range("H136").Value = 4520.3687
$ws.Range("I136").Value = 4494
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 13482
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -10932
$ws.Range("N136").Value = -20085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 483
$ws.Range("J80").Value = 454.26315
$ws.Range("L80").Value = 454.26315
$ws.Range("N80").Value = -2450.26315
$ws.Range("H83").Value = 483
$ws.Range("J83").Value = 454.26315
$ws.Range("L83").Value = 2271.31575
$ws.Range("N83").Value = -12255.31575
$ws.Range("H111").Value = 61999.5
$ws.Range("J111").Value = 61999.5
$ws.Range("L111").Value = 61999.5
$ws.Range("N111").Value = -70179.5
$ws.Range("H134").Value = 3271.7673
$ws.Range("I134").Value = 1425.6666
$ws.Range("J134").Value = 6387.0625
$ws.Range("K134").Value = 4276.9998
$ws.Range("L134").Value = 19161.1875
$ws.Range("M134").Value = -1741.9998
$ws.Range("N134").Value = -24231.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1666.6666
$ws.Range("I16").Value = 1462.4445
$ws.Range("K16").Value = 1462.4445
$ws.Range("M16").Value = -1175.4445
$ws.Range("H31").Value = 25261.842
$ws.Range("I31").Value = 2898.125
$ws.Range("J31").Value = 30231.555
$ws.Range("K31").Value = 2898.125
$ws.Range("L31").Value = 30231.555
$ws.Range("M31").Value = -2603.125
$ws.Range("N31").Value = -30821.555
$ws.Range("H34").Value = 25261.842
$ws.Range("I34").Value = 2898.125
$ws.Range("J34").Value = 30231.555
$ws.Range("K34").Value = 2898.125
$ws.Range("L34").Value = 30231.555
$ws.Range("M34").Value = -2696.125
$ws.Range("N34").Value = -30635.555
$ws.Range("H94").Value = 2058.6667
$ws.Range("I94").Value = 3333
$ws.Range("J94").Value = 1421.5
$ws.Range("K94").Value = 3333
$ws.Range("L94").Value = 1421.5
$ws.Range("M94").Value = -2882
$ws.Range("N94").Value = -2323.5
$ws.Range("H99").Value = 3504.1875
$ws.Range("I99").Value = 3006
$ws.Range("J99").Value = 4998.75
$ws.Range("K99").Value = 3006
$ws.Range("L99").Value = 4998.75
$ws.Range("M99").Value = -1508
$ws.Range("N99").Value = -7994.75
$ws.Range("H113").Value = 1666.6666
$ws.Range("I113").Value = 1462.4445
$ws.Range("K113").Value = 1462.4445
$ws.Range("M113").Value = 707.5554999999999
$ws.Range("H126").Value = 3504.1875
$ws.Range("I126").Value = 3006
$ws.Range("J126").Value = 4998.75
$ws.Range("K126").Value = 9018
$ws.Range("L126").Value = 14996.25
$ws.Range("M126").Value = -6548
$ws.Range("N126").Value = -19936.25
$ws.Range("H130").Value = 59240
$ws.Range("J130").Value = 59240
$ws.Range("L130").Value = 59240
$ws.Range("N130").Value = -69280
$ws.Range("H132").Value = 42499.98
$ws.Range("I132").Value = 27002.732
$ws.Range("J132").Value = 133269.58
$ws.Range("K132").Value = 81008.196
$ws.Range("L132").Value = 399808.74
$ws.Range("M132").Value = -78478.196
$ws.Range("N132").Value = -404868.74
$ws.Range("H141").Value = 673441.7
$ws.Range("J141").Value = 673441.7
$ws.Range("L141").Value = 673441.7
$ws.Range("N141").Value = -683801.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 48454.57
$ws.Range("I5").Value = 680.7692
$ws.Range("J5").Value = 126087
$ws.Range("K5").Value = 2042.3076
$ws.Range("L5").Value = 378261
$ws.Range("M5").Value = -1930.3076
$ws.Range("N5").Value = -378485
$ws.Range("H33").Value = 11190.223
$ws.Range("I33").Value = 76.666664
$ws.Range("K33").Value = 459.999984
$ws.Range("M33").Value = -176.999984
$ws.Range("H38").Value = 84.28570999999999
$ws.Range("J38").Value = 135
$ws.Range("L38").Value = 405
$ws.Range("N38").Value = -1099
$ws.Range("H107").Value = 990.875
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 1686.75
$ws.Range("K107").Value = 885
$ws.Range("L107").Value = 5060.25
$ws.Range("M107").Value = 1035
$ws.Range("N107").Value = -8900.25
$ws.Range("H135").Value = 48454.57
$ws.Range("I135").Value = 680.7692
$ws.Range("J135").Value = 126087
$ws.Range("K135").Value = 6126.922799999999
$ws.Range("L135").Value = 1134783
$ws.Range("M135").Value = -3591.922799999999
$ws.Range("N135").Value = -1139853

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 27030250
$ws.Range("I2").Value = 632.05
$ws.Range("K2").Value = 632.05
$ws.Range("M2").Value = -519.05
$ws.Range("H29").Value = 9841.6
$ws.Range("I29").Value = 8066.6665
$ws.Range("J29").Value = 12504
$ws.Range("K29").Value = 8066.6665
$ws.Range("L29").Value = 12504
$ws.Range("M29").Value = -7776.6665
$ws.Range("N29").Value = -13084
$ws.Range("H102").Value = 3777460.8
$ws.Range("I102").Value = 4831894.5
$ws.Range("J102").Value = 1911924.2
$ws.Range("K102").Value = 4831894.5
$ws.Range("L102").Value = 1911924.2
$ws.Range("M102").Value = -4830272.5
$ws.Range("N102").Value = -1915168.2
$ws.Range("H122").Value = 320071.5
$ws.Range("I122").Value = 470161.84
$ws.Range("K122").Value = 1410485.52
$ws.Range("M122").Value = -1408035.52
$ws.Range("H126").Value = 3484857.5
$ws.Range("I126").Value = 4547930
$ws.Range("J126").Value = 3091126.8
$ws.Range("K126").Value = 13643790
$ws.Range("L126").Value = 9273380.399999999
$ws.Range("M126").Value = -13641320
$ws.Range("N126").Value = -9278320.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7016.6665
$ws.Range("I40").Value = 4457.5713
$ws.Range("K40").Value = 4457.5713
$ws.Range("M40").Value = -4321.5713
$ws.Range("H55").Value = 1035.9131
$ws.Range("I55").Value = 1078.375
$ws.Range("J55").Value = 938.8570999999999
$ws.Range("K55").Value = 1078.375
$ws.Range("L55").Value = 938.8570999999999
$ws.Range("M55").Value = -905.375
$ws.Range("N55").Value = -1284.8571
$ws.Range("H57").Value = 41
$ws.Range("I57").Value = 41
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 41
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 525
$ws.Range("N57").ClearContents()
$ws.Range("H132").Value = 5727.2925
$ws.Range("I132").Value = 5757.472
$ws.Range("K132").Value = 17272.416
$ws.Range("M132").Value = -14742.416
$ws.Range("H136").Value = 31091.264
$ws.Range("I136").Value = 55716.26
$ws.Range("K136").Value = 167148.78
$ws.Range("M136").Value = -164598.78

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 23903.75
$ws.Range("J45").Value = 23903.75
$ws.Range("L45").Value = 23903.75
$ws.Range("N45").Value = -24885.75
$ws.Range("H100").Value = 1037.4615
$ws.Range("I100").Value = 1246.6666
$ws.Range("K100").Value = 2493.3332
$ws.Range("M100").Value = -1952.3332
$ws.Range("H132").Value = 18203152
$ws.Range("I132").Value = 20413540
$ws.Range("K132").Value = 61240620
$ws.Range("M132").Value = -61238090
$ws.Range("H136").Value = 3603.838
$ws.Range("I136").Value = 3203.3333
$ws.Range("J136").Value = 5320.2856
$ws.Range("K136").Value = 9609.999899999999
$ws.Range("L136").Value = 15960.8568
$ws.Range("M136").Value = -7059.999899999999
$ws.Range("N136").Value = -21060.8568
